$wb = $excel.ActiveWorkbook

$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f7ffd976057bad8ed1bc653b1d7145f477a2e28d/e2e/5ebd4d25-03b1-4ba0-a1c2-e488cc35cf20.md"
$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0eaba702b7ef3199f0b91ff6ebe2afb6ff9f7352/e2e/5ebd4d25-03b1-4ba0-a1c2-e488cc35cf20.md"
$errorMsg = "The version of handback file is not the latest, current: $currentUrl, latest: $latestUrl."

# zh-cn sheet: row 7 handback info just arrived
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Add($ws.Range("I7"), $latestUrl, "", "", "5ebd4d25-03b1-4ba0-a1c2-e488cc35cf20.md")
$ws.Range("J7").Value = "5ebd4d25-03b1-4ba0-a1c2-e488cc35cf20.0677ff1b4affe6fb753f2d37b0045ee5d48318b0.zh-cn.xlf"
$ws.Range("K7").Value = "2016-09-02 11:05:19"
$ws.Range("P7").Value = $errorMsg

# de-de sheet: row 7 handback info just arrived
$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Hyperlinks.Add($ws2.Range("I7"), $latestUrl, "", "", "5ebd4d25-03b1-4ba0-a1c2-e488cc35cf20.md")
$ws2.Range("J7").Value = "5ebd4d25-03b1-4ba0-a1c2-e488cc35cf20.0677ff1b4affe6fb753f2d37b0045ee5d48318b0.de-de.xlf"
$ws2.Range("K7").Value = "2016-09-02 11:05:27"
$ws2.Range("P7").Value = $errorMsg
